# Generate Report for Archive
#
# The localization status moved from "Ready for handoff" to "In Translation".
# That status string shows up on the Overview sheet (columns E/F, the per
# -language status for zh-cn / de-de) as well as on each language detail
# sheet (column C, "Status"). Updating the text makes the "Status" column
# narrower, since the new text is shorter than the old text, so the
# (auto-fitted) column width shrinks too.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: E2 (zh-cn status) and F2 (de-de status) ---
if ($overview.Range("E2").Value2 -eq $oldStatus) {
    $overview.Range("E2").Value = $newStatus
}
if ($overview.Range("F2").Value2 -eq $oldStatus) {
    $overview.Range("F2").Value = $newStatus
}

# --- zh-cn / de-de detail sheets: C2 ("Status" column) ---
if ($zhcn.Range("C2").Value2 -eq $oldStatus) {
    $zhcn.Range("C2").Value = $newStatus
}
if ($dede.Range("C2").Value2 -eq $oldStatus) {
    $dede.Range("C2").Value = $newStatus
}

# --- Shrink the now-narrower "Status" columns to match the shorter text ---
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth   # Overview!E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newWidth   # Overview!F (de-de)
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth = $newWidth        # de-de!C (Status)

Write-Output "Updated status text and column widths for archive report."
